$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.271338224411011
$ws.Range("B1").Value = 2.148607730865479
$ws.Range("C1").Value = 4.713882923126221
$ws.Range("D1").Value = 3.27408242225647
$ws.Range("E1").Value = 1.375074625015259
